# Add 2022-05-31 data: update nombre_aides (C) and montant_total (E) for
# specific rows in the Fonds de solidarite volet 1 regional/NAF dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=9;   C=69574;  E=191367649},
    @{Row=125; C=4597;   E=13136795},
    @{Row=150; C=95010;  E=278992606},
    @{Row=167; C=101529; E=194966146},
    @{Row=168; C=285008; E=1210390222},
    @{Row=169; C=562600; E=1284501384},
    @{Row=170; C=367381; E=2845626697},
    @{Row=171; C=115156; E=446646135},
    @{Row=173; C=54389;  E=151875013},
    @{Row=174; C=357235; E=1017861449},
    @{Row=175; C=125553; E=812943051},
    @{Row=177; C=96756;  E=174740928},
    @{Row=178; C=75360;  E=102747363},
    @{Row=179; C=235713; E=812627643},
    @{Row=180; C=141482; E=340835117},
    @{Row=181; C=7894;   E=11291524},
    @{Row=188; C=19707;  E=66147299},
    @{Row=204; C=4759;   E=11763349},
    @{Row=259; C=6240;   E=14751339}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
